$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Titile Text"
$ws.Range("B1").Value = "Title URL"

$ws.Range("A2").Value = "What Is Data Scraping? Definition & Usage"
$ws.Range("B2").Value = "https://www.okta.com/identity-101/data-scraping/"

$ws.Range("A3").Value = "What is data scraping?"
$ws.Range("B3").Value = "https://www.cloudflare.com/learning/bots/what-is-data-scraping/"

$ws.Range("A4").Value = "What is Web Scraping and How to Use It?"
$ws.Range("B4").Value = "https://www.geeksforgeeks.org/what-is-web-scraping-and-how-to-use-it/"
